# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Columns: D = Price (text, may look numeric), E = Volume(1h) % change (text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.853.66"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.293.00"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'102.75"
$ws.Range("E5").Value = "  +5.14%  "
$ws.Range("D6").Value = "'270.66"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "'46.32"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "'8.09"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "'15.55"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'0.857"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "2.286.97"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "43.742.11"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "'72.31"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'2.51"
$ws.Range("E21").Value = "  +10.58%  "
$ws.Range("D22").Value = "'233.47"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "'2.87"
$ws.Range("E23").Value = "  +13.15%  "
$ws.Range("D24").Value = "'9.28"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'11.25"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "'41.47"
$ws.Range("E27").Value = "  +7.97%  "
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'177.67"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").Value = "'21.81"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'0.0903"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "'5.50"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "'4.86"
$ws.Range("E34").Value = "  +9.51%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "'3.57"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "'0.237"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").Value = "'2.33"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "'65.93"
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("D43").Value = "'12.24"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'5.27"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "'8.81"
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "'1.24"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").Value = "'99.09"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "'0.443"
$ws.Range("E49").Value = "  +6.07%  "
$ws.Range("E50").Value = "  +11.04%  "
$ws.Range("D51").Value = "2.515.10"
$ws.Range("E51").Value = "  -1.11%  "
